$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "se aceptan los sets": the FOLLOW sets in column B are recomputed/renumbered.
# Set every affected cell to its new literal text; Excel's shared-string table
# will automatically de-duplicate identical values (e.g. the many cells that
# now read "[88]").
$ws.Range("B4").Value = "[88]"
$ws.Range("B11").Value = "[88]"
$ws.Range("B12").Value = "[12, 13, 15]"
$ws.Range("B13").Value = "[12, 13, 15]"
$ws.Range("B14").Value = "[12, 13, 15]"
$ws.Range("B16").Value = "[88]"
$ws.Range("B18").Value = "[88]"
$ws.Range("B20").Value = "[20]"
$ws.Range("B21").Value = "[88]"
$ws.Range("B23").Value = "[88]"
$ws.Range("B25").Value = "[88]"
$ws.Range("B27").Value = "[27]"
$ws.Range("B28").Value = "[88]"
$ws.Range("B30").Value = "[30]"
$ws.Range("B31").Value = "[88]"
$ws.Range("B33").Value = "[88]"
$ws.Range("B35").Value = "[88]"
$ws.Range("B37").Value = "[37]"
$ws.Range("B38").Value = "[88]"
$ws.Range("B40").Value = "[88]"
$ws.Range("B42").Value = "[42]"
$ws.Range("B44").Value = "[88]"
$ws.Range("B46").Value = "[46]"
$ws.Range("B48").Value = "[88]"
$ws.Range("B50").Value = "[50]"
$ws.Range("B52").Value = "[88]"
$ws.Range("B54").Value = "[54]"
$ws.Range("B56").Value = "[88]"
$ws.Range("B58").Value = "[58]"
$ws.Range("B59").Value = "[88]"
$ws.Range("B61").Value = "[61]"
$ws.Range("B62").Value = "[88]"
$ws.Range("B64").Value = "[88]"
$ws.Range("B66").Value = "[88]"
$ws.Range("B68").Value = "[88]"
$ws.Range("B70").Value = "[88]"
$ws.Range("B72").Value = "[88]"
$ws.Range("B74").Value = "[88]"
$ws.Range("B76").Value = "[88]"
$ws.Range("B78").Value = "[88]"
$ws.Range("B80").Value = "[80]"
$ws.Range("B81").Value = "[88]"
$ws.Range("B83").Value = "[88]"
$ws.Range("B85").Value = "[88]"
$ws.Range("B87").Value = "[87]"

# The last two rows (89 and 90) are no longer needed; remove them entirely
# so the data range shrinks from A1:B90 to A1:B88.
$ws.Rows.Item(90).Delete()
$ws.Rows.Item(89).Delete()
